# Applies the weekly Fruit/Vegetable price update described in the commit:
# "Fruta / hortaliza, semanal" -- updates rows 202-223 with revised values,
# replaces row 224, and appends new rows 225-227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 202-223 (only the cells that actually changed) ---
# Row 202
$ws.Cells.Item(202, 4).Value = 44578
$ws.Cells.Item(202, 13).Value = 70
$ws.Cells.Item(202, 14).Value = 6500
$ws.Cells.Item(202, 15).Value = 6500
$ws.Cells.Item(202, 16).Value = 6500
$ws.Cells.Item(202, 19).Value = 929

# Row 203
$ws.Cells.Item(203, 4).Value = 44578
$ws.Cells.Item(203, 12).Value = "Primera"
$ws.Cells.Item(203, 13).Value = 80

# Row 204
$ws.Cells.Item(204, 4).Value = 44578
$ws.Cells.Item(204, 12).Value = "Segunda"
$ws.Cells.Item(204, 13).Value = 85
$ws.Cells.Item(204, 14).Value = 4000
$ws.Cells.Item(204, 15).Value = 4000
$ws.Cells.Item(204, 16).Value = 4000
$ws.Cells.Item(204, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(204, 19).Value = 571

# Row 205
$ws.Cells.Item(205, 4).Value = 44490
$ws.Cells.Item(205, 12).Value = "Especial"
$ws.Cells.Item(205, 13).Value = 98
$ws.Cells.Item(205, 14).Value = 7000
$ws.Cells.Item(205, 15).Value = 7000
$ws.Cells.Item(205, 16).Value = 7000
$ws.Cells.Item(205, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(205, 19).Value = 1000

# Row 206
$ws.Cells.Item(206, 4).Value = 44490
$ws.Cells.Item(206, 12).Value = "Segunda"
$ws.Cells.Item(206, 13).Value = 87
$ws.Cells.Item(206, 14).Value = 5000
$ws.Cells.Item(206, 15).Value = 5000
$ws.Cells.Item(206, 16).Value = 5000
$ws.Cells.Item(206, 19).Value = 714

# Row 207
$ws.Cells.Item(207, 4).Value = 44491
$ws.Cells.Item(207, 12).Value = "Especial"
$ws.Cells.Item(207, 13).Value = 70
$ws.Cells.Item(207, 14).Value = 7000
$ws.Cells.Item(207, 15).Value = 7000
$ws.Cells.Item(207, 16).Value = 7000
$ws.Cells.Item(207, 18).Value = "Provincia de Cautín"
$ws.Cells.Item(207, 19).Value = 1000

# Row 208
$ws.Cells.Item(208, 4).Value = 44491
$ws.Cells.Item(208, 12).Value = "Segunda"
$ws.Cells.Item(208, 13).Value = 60
$ws.Cells.Item(208, 14).Value = 5000
$ws.Cells.Item(208, 15).Value = 5000
$ws.Cells.Item(208, 16).Value = 5000
$ws.Cells.Item(208, 18).Value = "Provincia de Cautín"
$ws.Cells.Item(208, 19).Value = 714

# Row 209
$ws.Cells.Item(209, 4).Value = 44293
$ws.Cells.Item(209, 13).Value = 50
$ws.Cells.Item(209, 14).Value = 10000
$ws.Cells.Item(209, 15).Value = 10000
$ws.Cells.Item(209, 16).Value = 10000
$ws.Cells.Item(209, 19).Value = 1429

# Row 210
$ws.Cells.Item(210, 4).Value = 44293
$ws.Cells.Item(210, 12).Value = "Primera"
$ws.Cells.Item(210, 13).Value = 50
$ws.Cells.Item(210, 14).Value = 9000
$ws.Cells.Item(210, 15).Value = 9000
$ws.Cells.Item(210, 16).Value = 9000
$ws.Cells.Item(210, 19).Value = 1286

# Row 211
$ws.Cells.Item(211, 4).Value = 44266
$ws.Cells.Item(211, 13).Value = 78
$ws.Cells.Item(211, 14).Value = 8500
$ws.Cells.Item(211, 15).Value = 8500
$ws.Cells.Item(211, 16).Value = 8500
$ws.Cells.Item(211, 19).Value = 1214

# Row 212
$ws.Cells.Item(212, 4).Value = 44533
$ws.Cells.Item(212, 13).Value = 139
$ws.Cells.Item(212, 14).Value = 5500
$ws.Cells.Item(212, 15).Value = 6000
$ws.Cells.Item(212, 16).Value = 5806
$ws.Cells.Item(212, 19).Value = 829

# Row 213
$ws.Cells.Item(213, 4).Value = 44533
$ws.Cells.Item(213, 13).Value = 78
$ws.Cells.Item(213, 14).Value = 4000
$ws.Cells.Item(213, 15).Value = 4000
$ws.Cells.Item(213, 16).Value = 4000
$ws.Cells.Item(213, 19).Value = 571

# Row 214
$ws.Cells.Item(214, 4).Value = 44264
$ws.Cells.Item(214, 13).Value = 50

# Row 215
$ws.Cells.Item(215, 4).Value = 44494
$ws.Cells.Item(215, 13).Value = 85
$ws.Cells.Item(215, 14).Value = 7000
$ws.Cells.Item(215, 15).Value = 7000
$ws.Cells.Item(215, 16).Value = 7000
$ws.Cells.Item(215, 19).Value = 1000

# Row 216
$ws.Cells.Item(216, 4).Value = 44494
$ws.Cells.Item(216, 13).Value = 47
$ws.Cells.Item(216, 14).Value = 5000
$ws.Cells.Item(216, 15).Value = 5000
$ws.Cells.Item(216, 16).Value = 5000
$ws.Cells.Item(216, 19).Value = 714

# Row 217
$ws.Cells.Item(217, 4).Value = 44279
$ws.Cells.Item(217, 13).Value = 70
$ws.Cells.Item(217, 14).Value = 8000
$ws.Cells.Item(217, 15).Value = 8000
$ws.Cells.Item(217, 16).Value = 8000
$ws.Cells.Item(217, 19).Value = 1143

# Row 218
$ws.Cells.Item(218, 4).Value = 44525
$ws.Cells.Item(218, 13).Value = 87

# Row 219
$ws.Cells.Item(219, 4).Value = 44525
$ws.Cells.Item(219, 13).Value = 89

# Row 220
$ws.Cells.Item(220, 4).Value = 44354
$ws.Cells.Item(220, 12).Value = "Especial"
$ws.Cells.Item(220, 13).Value = 45
$ws.Cells.Item(220, 14).Value = 14000
$ws.Cells.Item(220, 15).Value = 14000
$ws.Cells.Item(220, 16).Value = 14000
$ws.Cells.Item(220, 19).Value = 2000

# Row 221
$ws.Cells.Item(221, 4).Value = 44503
$ws.Cells.Item(221, 13).Value = 90

# Row 222
$ws.Cells.Item(222, 4).Value = 44503
$ws.Cells.Item(222, 13).Value = 70

# Row 223
$ws.Cells.Item(223, 4).Value = 44462
$ws.Cells.Item(223, 12).Value = "Primera"
$ws.Cells.Item(223, 13).Value = 50
$ws.Cells.Item(223, 14).Value = 17000
$ws.Cells.Item(223, 15).Value = 17000
$ws.Cells.Item(223, 16).Value = 17000
$ws.Cells.Item(223, 19).Value = 2429

# --- Row 224: content fully replaced with a new record ---
$ws.Cells.Item(224, 1).Value = 3
$ws.Cells.Item(224, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 44512
$ws.Cells.Item(224, 5).Value = 5
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100101
$ws.Cells.Item(224, 8).Value = "Berries"
$ws.Cells.Item(224, 9).Value = 100112025
$ws.Cells.Item(224, 10).Value = "Frutilla"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Especial"
$ws.Cells.Item(224, 13).Value = 85
$ws.Cells.Item(224, 14).Value = 6000
$ws.Cells.Item(224, 15).Value = 6000
$ws.Cells.Item(224, 16).Value = 6000
$ws.Cells.Item(224, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(224, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(224, 19).Value = 857
$ws.Cells.Item(224, 20).Value = 7

# --- New rows 225-227 appended at the end ---
# Row 225
$ws.Cells.Item(225, 1).Value = 3
$ws.Cells.Item(225, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(225, 3).Value = "Coquimbo"
$ws.Cells.Item(225, 4).Value = 44512
$ws.Cells.Item(225, 5).Value = 5
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100101
$ws.Cells.Item(225, 8).Value = "Berries"
$ws.Cells.Item(225, 9).Value = 100112025
$ws.Cells.Item(225, 10).Value = "Frutilla"
$ws.Cells.Item(225, 11).Value = "Sin especificar"
$ws.Cells.Item(225, 12).Value = "Segunda"
$ws.Cells.Item(225, 13).Value = 80
$ws.Cells.Item(225, 14).Value = 4000
$ws.Cells.Item(225, 15).Value = 4000
$ws.Cells.Item(225, 16).Value = 4000
$ws.Cells.Item(225, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(225, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(225, 19).Value = 571
$ws.Cells.Item(225, 20).Value = 7
$ws.Cells.Item(225, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 226
$ws.Cells.Item(226, 1).Value = 3
$ws.Cells.Item(226, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 44511
$ws.Cells.Item(226, 5).Value = 5
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100101
$ws.Cells.Item(226, 8).Value = "Berries"
$ws.Cells.Item(226, 9).Value = 100112025
$ws.Cells.Item(226, 10).Value = "Frutilla"
$ws.Cells.Item(226, 11).Value = "Sin especificar"
$ws.Cells.Item(226, 12).Value = "Especial"
$ws.Cells.Item(226, 13).Value = 125
$ws.Cells.Item(226, 14).Value = 6000
$ws.Cells.Item(226, 15).Value = 6000
$ws.Cells.Item(226, 16).Value = 6000
$ws.Cells.Item(226, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(226, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(226, 19).Value = 857
$ws.Cells.Item(226, 20).Value = 7
$ws.Cells.Item(226, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 227
$ws.Cells.Item(227, 1).Value = 3
$ws.Cells.Item(227, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(227, 3).Value = "Coquimbo"
$ws.Cells.Item(227, 4).Value = 44511
$ws.Cells.Item(227, 5).Value = 5
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100101
$ws.Cells.Item(227, 8).Value = "Berries"
$ws.Cells.Item(227, 9).Value = 100112025
$ws.Cells.Item(227, 10).Value = "Frutilla"
$ws.Cells.Item(227, 11).Value = "Sin especificar"
$ws.Cells.Item(227, 12).Value = "Segunda"
$ws.Cells.Item(227, 13).Value = 80
$ws.Cells.Item(227, 14).Value = 4000
$ws.Cells.Item(227, 15).Value = 4000
$ws.Cells.Item(227, 16).Value = 4000
$ws.Cells.Item(227, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(227, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(227, 19).Value = 571
$ws.Cells.Item(227, 20).Value = 7
$ws.Cells.Item(227, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

